$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.513.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.889.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.77%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.503"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.887.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.52%  "

$ws.Range("E11").Value = "  -2.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.428"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.69%  "

$ws.Range("E13").Value = "  -1.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.27%  "

$ws.Range("E15").Value = "  -0.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.367.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.491.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.883.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "430.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.652"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.88%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.08%  "

$ws.Range("E25").Value = "  +0.43%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -12.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000105"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.77%  "

$ws.Range("E31").Value = "  -4.47%  "

$ws.Range("E32").Value = "  -9.18%  "

$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("E34").Value = "  -2.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.954"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.63%  "

$ws.Range("E40").Value = "  -9.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.69%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.113"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.266"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.98%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.685.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "132.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0333"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.40%  "

$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "341.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.07%  "

$ws.Range("E50").Value = "  -2.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.64%  "
